$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Must have properties")

$ws.Range("A19").Value = "file.file_size"
$ws.Range("B19").Value = "file"
$ws.Range("C19").Value = "file_size"

$ws.Range("A20").Value = "file.md5sum"
$ws.Range("B20").Value = "file"
$ws.Range("C20").Value = "md5sum"
